$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Julio de 2020 a las 07:58"

# --- Row 6: India (casos totales etc. refreshed) ---
$ws.Range("B6").Value = 1533936
$ws.Range("C6").Value = 1801
$ws.Range("D6").Value = 989624
$ws.Range("E6").Value = 510072
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 16
$ws.Range("H6").Value = 34240

# --- Row 15: Pakistan (refreshed) ---
$ws.Range("B15").Value = 276288
$ws.Range("C15").Value = 1063
$ws.Range("D15").Value = 244883
$ws.Range("E15").Value = 25513
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 27
$ws.Range("H15").Value = 5892

# --- Rows 56 & 57: Kirguistan overtakes Ghana in total cases, so they swap order ---
# Row 56 becomes Kirguistan with its refreshed numbers
$ws.Range("A56").Value = "Kirguistan"
$ws.Range("B56").Value = 34592
$ws.Range("C56").Value = 748
$ws.Range("D56").Value = 22296
$ws.Range("E56").Value = 10949
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 18
$ws.Range("H56").Value = 1347

# Row 57 becomes Ghana, keeping its previous (unchanged) numbers
$ws.Range("A57").Value = "Ghana"
$ws.Range("B57").Value = 34406
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 30621
$ws.Range("E57").Value = 3617
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 168

# --- Row 64: Uzbekistan (refreshed) ---
$ws.Range("B64").Value = 22169
$ws.Range("C64").Value = 276
$ws.Range("D64").Value = 12265
$ws.Range("E64").Value = 9777
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 127

# --- Row 91: Haiti (refreshed) ---
$ws.Range("B91").Value = 7371
$ws.Range("C91").Value = 31
$ws.Range("D91").Value = 4467
$ws.Range("E91").Value = 2746

# --- Row 109: Tailandia (refreshed) ---
$ws.Range("B109").Value = 3298
$ws.Range("C109").Value = 1
$ws.Range("E109").Value = 129
